$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "22.386.49"
Set-TextValue "E2" "  -4.27%  "
Set-TextValue "D3" "1.569.90"
Set-TextValue "E3" "  -4.07%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "E5" "  -0.11%  "
Set-TextValue "D6" "289.47"
Set-TextValue "E6" "  -3.36%  "
Set-TextValue "D7" "0.3690"
Set-TextValue "E7" "  -2.41%  "
Set-TextValue "D8" "49.36"
Set-TextValue "E8" "  -2.00%  "
Set-TextValue "D9" "0.3377"
Set-TextValue "E9" "  -4.11%  "
Set-TextValue "D10" "1.165"
Set-TextValue "E10" "  -3.52%  "
Set-TextValue "D11" "0.07604"
Set-TextValue "E11" "  -5.54%  "
Set-TextValue "D12" "0.9999"
Set-TextValue "E12" "  -0.18%  "
Set-TextValue "D13" "21.19"
Set-TextValue "E13" "  -3.42%  "
Set-TextValue "D14" "6.055"
Set-TextValue "E14" "  -4.42%  "
Set-TextValue "D15" "6.895"
Set-TextValue "E15" "  -4.96%  "
Set-TextValue "D16" "1.572.75"
Set-TextValue "E16" "  -4.05%  "
Set-TextValue "D17" "0.00001134"
Set-TextValue "E17" "  -5.66%  "
Set-TextValue "D18" "89.00"
Set-TextValue "E18" "  -7.24%  "
Set-TextValue "D19" "0.06754"
Set-TextValue "E19" "  -2.86%  "
Set-TextValue "E20" "  -0.04%  "
Set-TextValue "D21" "6.238"
Set-TextValue "E21" "  -6.94%  "
Set-TextValue "B22" "Avalanche"
Set-TextValue "C22" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D22" "16.54"
Set-TextValue "E22" "  -4.45%  "
Set-TextValue "B23" "BitDAO"
Set-TextValue "C23" "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
Set-TextValue "D23" "0.5321"
Set-TextValue "E23" "  -6.91%  "
Set-TextValue "D24" "11.96"
Set-TextValue "E24" "  -2.71%  "
Set-TextValue "D25" "22.408.95"
Set-TextValue "E25" "  -4.24%  "
Set-TextValue "D26" "2.382"
Set-TextValue "E26" "  -3.70%  "
Set-TextValue "D27" "2.987"
Set-TextValue "E27" "  +3.42%  "
Set-TextValue "D28" "19.92"
Set-TextValue "E28" "  -4.24%  "
Set-TextValue "D29" "145.07"
Set-TextValue "E29" "  -5.01%  "
Set-TextValue "D30" "4.967"
Set-TextValue "E30" "  -4.24%  "
Set-TextValue "D31" "125.23"
Set-TextValue "E31" "  -5.66%  "
Set-TextValue "D32" "1.746.02"
Set-TextValue "E32" "  -4.08%  "
Set-TextValue "D33" "1.039"
Set-TextValue "E33" "  +6.77%  "
Set-TextValue "D34" "6.258"
Set-TextValue "E34" "  -8.31%  "
Set-TextValue "E35" "  -6.17%  "
Set-TextValue "D36" "10.31"
Set-TextValue "E36" "  -9.02%  "
Set-TextValue "D37" "0.08458"
Set-TextValue "E37" "  -3.07%  "
Set-TextValue "D38" "0.02534"
Set-TextValue "E38" "  -6.23%  "
Set-TextValue "D39" "0.2330"
Set-TextValue "E39" "  -3.70%  "
Set-TextValue "D40" "5.550"
Set-TextValue "E40" "  -5.55%  "
Set-TextValue "D41" "0.06509"
Set-TextValue "E41" "  -4.20%  "
Set-TextValue "D42" "11.82"
Set-TextValue "E42" "  -9.05%  "
Set-TextValue "D43" "1.245"
Set-TextValue "E43" "  -4.43%  "
Set-TextValue "D44" "0.6366"
Set-TextValue "E44" "  -6.76%  "
Set-TextValue "D45" "14.34"
Set-TextValue "E45" "  -7.22%  "
Set-TextValue "D46" "0.9998"
Set-TextValue "E46" "  -0.13%  "
Set-TextValue "D47" "0.5983"
Set-TextValue "E47" "  -5.30%  "
Set-TextValue "D48" "3.754"
Set-TextValue "E48" "  -3.77%  "
Set-TextValue "D49" "2.130"
Set-TextValue "E49" "  -4.96%  "
Set-TextValue "D50" "1.259"
Set-TextValue "E50" "  +3.16%  "
Set-TextValue "D51" "123.42"
Set-TextValue "E51" "  -2.61%  "
